$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.781.01"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "3.096.74"
$ws.Range("E3").Value = "  +4.93%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'579.20"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").Value = "'172.47"
$ws.Range("E6").Value = "  +5.28%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.091.98"
$ws.Range("E8").Value = "  +4.98%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("D12").Value = "'0.481"
$ws.Range("E12").Value = "  +3.57%  "
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "'37.36"
$ws.Range("E14").Value = "  +6.56%  "
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "3.607.17"
$ws.Range("E16").Value = "  +5.04%  "
$ws.Range("D17").Value = "66.734.05"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "'7.19"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "3.096.56"
$ws.Range("E19").Value = "  +5.21%  "
$ws.Range("D20").Value = "'16.29"
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("D21").Value = "'480.81"
$ws.Range("E21").Value = "  +7.74%  "
$ws.Range("D22").Value = "'0.715"
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("D23").Value = "'7.54"
$ws.Range("E23").Value = "  +3.59%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'84.01"
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'13.22"
$ws.Range("E25").Value = "  +7.31%  "
$ws.Range("E26").Value = "  +5.36%  "
$ws.Range("D27").Value = "'10.03"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'8.00"
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("D32").Value = "'28.81"
$ws.Range("E32").Value = "  +5.43%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "'0.115"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  +2.95%  "
$ws.Range("D37").Value = "'0.990"
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("D38").Value = "'47.92"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("E39").Value = "  +6.79%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.318"
$ws.Range("E40").Value = "  +5.06%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'50.01"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").Value = "'8.66"
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("D45").Value = "2.841.39"
$ws.Range("E45").Value = "  +6.20%  "
$ws.Range("D46").Value = "'0.0360"
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("D47").Value = "'383.62"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "'135.32"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D50").Value = "'25.04"
$ws.Range("E50").Value = "  +4.58%  "
$ws.Range("E51").Value = "  +2.65%  "
